$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$formula = "selected(`${(prefixo)_DGE_SQE_B1_P13_1_classes_estapas}, 'classe_1') or selected(`${(prefixo)_DGE_SQE_B1_P13_1_classes_estapas}, 'classe_2') or selected(`${(prefixo)_DGE_SQE_B1_P13_1_classes_estapas}, 'classe_3')"

# Row 115: fix the relevant-condition formula, then rename the variable
$ws.Range("B115").Value = $formula
$ws.Range("A115").Value = "pessoal_direcao_escola_subdirector_director_tecnico"

# Row 116: fix the relevant-condition formula, then rename the variable
$ws.Range("B116").Value = $formula
$ws.Range("A116").Value = "pessoal_direcao_escola_coordenador_pedagógico"

# Match the styling used for the rest of the "relevante" column (B) in this block
$ws.Range("B115").Style = $ws.Range("B117").Style
$ws.Range("B116").Style = $ws.Range("B117").Style

# A115/A116 now use the plain default style, matching the other freshly-typed variable cells
$ws.Range("A115").Style = $ws.Range("A113").Style
$ws.Range("A116").Style = $ws.Range("A113").Style

$ws.Range("A116").Select()
$ws.Application.ActiveWindow.ScrollRow = 100
